$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12 and 13 swapped places in the ranking (Wrapped liquid staked Ether 2.0
# moved above Wrapped Ether)
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# Refreshed Price / Volume(1h) figures. A few Price values look like plain numbers
# (e.g. "1.00", "3.30") but the source data is plain text, so a leading apostrophe
# keeps Excel from re-parsing them as numeric and dropping the trailing zeros; the
# style is reset straight back to Normal so only the cell VALUE changes.
$ws.Range("D2").Value = "26.021.31"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "1.643.05"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'216.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "1.871.28"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "1.670.50"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "'63.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "26.109.72"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'195.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").Value = "'1.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "'143.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "'15.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").Value = "'0.0498"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'3.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").Value = "1.132.91"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("D42").Value = "'99.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'0.796"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "1.780.57"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "0.0₆0117"
$ws.Range("E45").Value = "  +6.92%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").Value = "'0.0524"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'7.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -0.14%  "
